# TC06 - Job Categories feature
$wb = $excel.ActiveWorkbook

# Rename existing sheets to include numeric prefixes
$wsAccount = $wb.Worksheets.Item("Add New Account")
$wsAccount.Name = "01 - Add New Account"

$wsJobTitle = $wb.Worksheets.Item("Add Job Title")
$wsJobTitle.Name = "02 - Add Job Title"

# Add new worksheet for Job Category after the Job Title sheet (i.e. at the end)
$wsJobCategory = $wb.Worksheets.Add($null, $wsJobTitle)
$wsJobCategory.Name = "04 - Add Job Category"

# Populate the new sheet with header + data
$wsJobCategory.Range("A1").Value = "jobCategoryName"
$wsJobCategory.Range("A2").Value = "IT Professional"
$wsJobCategory.Range("A3").Value = "Education"
$wsJobCategory.Range("A4").Value = "Health care"
$wsJobCategory.Range("A5").Value = "Painter"
$wsJobCategory.Range("A6").Value = "Barber"

# Set the new sheet as the active/selected tab
$wsJobCategory.Select()
